$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 306.875
$ws.Range("I6").Value = 308.57144
$ws.Range("J6").Value = 295
$ws.Range("K6").Value = 925.71432
$ws.Range("L6").Value = 885
$ws.Range("M6").Value = -813.71432
$ws.Range("N6").Value = -1109
$ws.Range("H138").Value = 5150.851
$ws.Range("I138").Value = 2611.625
$ws.Range("J138").Value = 5671.718
$ws.Range("K138").Value = 7834.875
$ws.Range("L138").Value = 17015.154
$ws.Range("M138").Value = -2694.875
$ws.Range("N138").Value = -27295.154

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1479.65
$ws.Range("I32").Value = 1479.65
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1479.65
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1192.65
$ws.Range("H61").Value = 5000
$ws.Range("I61").Value = 5000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 5000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4788
$ws.Range("N61").ClearContents()
$ws.Range("H122").Value = 4106.65
$ws.Range("I122").Value = 3743.9092
$ws.Range("J122").Value = 4550
$ws.Range("K122").Value = 11231.7276
$ws.Range("L122").Value = 13650
$ws.Range("M122").Value = -8781.7276
$ws.Range("N122").Value = -18550
$ws.Range("H132").Value = 4653.0435
$ws.Range("I132").Value = 4801.05
$ws.Range("J132").Value = 3666.3333
$ws.Range("K132").Value = 14403.15
$ws.Range("L132").Value = 10998.9999
$ws.Range("M132").Value = -11873.15
$ws.Range("N132").Value = -16058.9999
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 15000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -12450
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 741337.9399999999
$ws.Range("I86").Value = 1064486
$ws.Range("J86").Value = 2714
$ws.Range("K86").Value = 1064486
$ws.Range("L86").Value = 2714
$ws.Range("M86").Value = -1063363
$ws.Range("N86").Value = -4960
$ws.Range("H89").Value = 741337.9399999999
$ws.Range("I89").Value = 1064486
$ws.Range("J89").Value = 2714
$ws.Range("K89").Value = 5322430
$ws.Range("L89").Value = 13570
$ws.Range("M89").Value = -5316814
$ws.Range("N89").Value = -24802
$ws.Range("H134").Value = 42915.52
$ws.Range("I134").Value = 5639.9565
$ws.Range("J134").Value = 257250
$ws.Range("K134").Value = 16919.8695
$ws.Range("L134").Value = 771750
$ws.Range("M134").Value = -14384.8695
$ws.Range("N134").Value = -776820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8437.071
$ws.Range("I16").Value = 2510.6
$ws.Range("J16").Value = 23253.25
$ws.Range("K16").Value = 2510.6
$ws.Range("L16").Value = 23253.25
$ws.Range("M16").Value = -2223.6
$ws.Range("N16").Value = -23827.25
$ws.Range("H31").Value = 48675.043
$ws.Range("I31").Value = 1702.3
$ws.Range("J31").Value = 84807.92
$ws.Range("K31").Value = 1702.3
$ws.Range("L31").Value = 84807.92
$ws.Range("M31").Value = -1407.3
$ws.Range("N31").Value = -85397.92
$ws.Range("H34").Value = 48675.043
$ws.Range("I34").Value = 1702.3
$ws.Range("J34").Value = 84807.92
$ws.Range("K34").Value = 1702.3
$ws.Range("L34").Value = 84807.92
$ws.Range("M34").Value = -1500.3
$ws.Range("N34").Value = -85211.92
$ws.Range("H39").Value = 12000
$ws.Range("I39").Value = 12000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 12000
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -11609
$ws.Range("H49").Value = 12000
$ws.Range("I49").Value = 12000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 12000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -11818
$ws.Range("H58").Value = 3901.8333
$ws.Range("I58").Value = 4803.6665
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 4803.6665
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -4600.6665
$ws.Range("N58").Value = -3406
$ws.Range("H107").Value = 529.7241
$ws.Range("I107").Value = 492.21738
$ws.Range("J107").Value = 673.5
$ws.Range("K107").Value = 492.21738
$ws.Range("L107").Value = 673.5
$ws.Range("M107").Value = 1427.78262
$ws.Range("N107").Value = -4513.5
$ws.Range("H113").Value = 8437.071
$ws.Range("I113").Value = 2510.6
$ws.Range("J113").Value = 23253.25
$ws.Range("K113").Value = 2510.6
$ws.Range("L113").Value = 23253.25
$ws.Range("M113").Value = -340.5999999999999
$ws.Range("N113").Value = -27593.25
$ws.Range("H132").Value = 2924.2307
$ws.Range("I132").Value = 2776.6667
$ws.Range("J132").Value = 3256.25
$ws.Range("K132").Value = 8330.000100000001
$ws.Range("L132").Value = 9768.75
$ws.Range("M132").Value = -5800.000100000001
$ws.Range("N132").Value = -14828.75
$ws.Range("H134").Value = 773781.6
$ws.Range("I134").Value = 4917.4546
$ws.Range("J134").Value = 5002534.5
$ws.Range("K134").Value = 14752.3638
$ws.Range("L134").Value = 15007603.5
$ws.Range("M134").Value = -12217.3638
$ws.Range("N134").Value = -15012673.5
$ws.Range("H136").Value = 3901.8333
$ws.Range("I136").Value = 4803.6665
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 14410.9995
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -11860.9995
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 181
$ws.Range("I26").Value = 181
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 543
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -255
$ws.Range("H50").Value = 495.42856
$ws.Range("I50").Value = 590
$ws.Range("J50").Value = 259
$ws.Range("K50").Value = 1770
$ws.Range("L50").Value = 777
$ws.Range("M50").Value = -1289
$ws.Range("N50").Value = -1739
$ws.Range("H53").Value = 495.42856
$ws.Range("I53").Value = 590
$ws.Range("J53").Value = 259
$ws.Range("K53").Value = 1770
$ws.Range("L53").Value = 777
$ws.Range("M53").Value = -1289
$ws.Range("N53").Value = -1739
$ws.Range("H80").Value = 1714.7142
$ws.Range("I80").Value = 1501
$ws.Range("J80").Value = 1800.2
$ws.Range("K80").Value = 4503
$ws.Range("L80").Value = 5400.6
$ws.Range("M80").Value = -3567
$ws.Range("N80").Value = -7272.6
$ws.Range("H83").Value = 1714.7142
$ws.Range("I83").Value = 1501
$ws.Range("J83").Value = 1800.2
$ws.Range("K83").Value = 13509
$ws.Range("L83").Value = 16201.8
$ws.Range("M83").Value = -8829
$ws.Range("N83").Value = -25561.8
$ws.Range("H96").Value = 668665
$ws.Range("I96").Value = 999995
$ws.Range("J96").Value = 503000
$ws.Range("K96").Value = 2999985
$ws.Range("L96").Value = 1509000
$ws.Range("M96").Value = -2997926
$ws.Range("N96").Value = -1513118
$ws.Range("H139").Value = 8071.357
$ws.Range("I139").Value = 8999
$ws.Range("J139").Value = 8000
$ws.Range("K139").Value = 26997
$ws.Range("L139").Value = 24000
$ws.Range("M139").Value = -21857
$ws.Range("N139").Value = -34280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 579.4375
$ws.Range("I97").Value = 579.4375
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 579.4375
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -83.4375
$ws.Range("H102").Value = 2067.2593
$ws.Range("I102").Value = 1424.619
$ws.Range("J102").Value = 4316.5
$ws.Range("K102").Value = 1424.619
$ws.Range("L102").Value = 4316.5
$ws.Range("M102").Value = 197.3810000000001
$ws.Range("N102").Value = -7560.5
$ws.Range("H113").Value = 591104.3
$ws.Range("I113").Value = 1001727.3
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 1001727.3
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -999557.3
$ws.Range("N113").Value = -8840
$ws.Range("H132").Value = 108938.2
$ws.Range("I132").Value = 10547.75
$ws.Range("J132").Value = 502500
$ws.Range("K132").Value = 31643.25
$ws.Range("L132").Value = 1507500
$ws.Range("M132").Value = -29113.25
$ws.Range("N132").Value = -1512560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 199000
$ws.Range("I23").Value = 199000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 199000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -198770
$ws.Range("H61").Value = 3789.6667
$ws.Range("I61").Value = 3789.6667
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3789.6667
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3587.6667
$ws.Range("H113").Value = 3789.6667
$ws.Range("I113").Value = 3789.6667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3789.6667
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1619.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 12122.6
$ws.Range("I74").Value = 6998
$ws.Range("J74").Value = 13403.75
$ws.Range("K74").Value = 6998
$ws.Range("L74").Value = 13403.75
$ws.Range("M74").Value = -6062
$ws.Range("N74").Value = -15275.75
$ws.Range("H77").Value = 12122.6
$ws.Range("I77").Value = 6998
$ws.Range("J77").Value = 13403.75
$ws.Range("K77").Value = 20994
$ws.Range("L77").Value = 40211.25
$ws.Range("M77").Value = -16314
$ws.Range("N77").Value = -49571.25
$ws.Range("H107").Value = 1130.0312
$ws.Range("I107").Value = 1381.3914
$ws.Range("J107").Value = 487.66666
$ws.Range("K107").Value = 4144.174199999999
$ws.Range("L107").Value = 1462.99998
$ws.Range("M107").Value = -2224.174199999999
$ws.Range("N107").Value = -5302.999980000001
$ws.Range("H113").Value = 300
$ws.Range("I113").Value = 300.25
$ws.Range("J113").Value = 299
$ws.Range("K113").Value = 900.75
$ws.Range("L113").Value = 897
$ws.Range("M113").Value = 1269.25
$ws.Range("N113").Value = -5237
$ws.Range("H136").Value = 175329.92
$ws.Range("I136").Value = 40738.273
$ws.Range("J136").Value = 668832.7
$ws.Range("K136").Value = 122214.819
$ws.Range("L136").Value = 2006498.1
$ws.Range("M136").Value = -119664.819
$ws.Range("N136").Value = -2011598.1

Write-Host "Applied all Jenova_Profits updates"